$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: the maturity note in column B is replaced with a longer note,
#     and the row grows to a taller wrapped height. ---
$ws.Cells.Item(2,2).Value = 'From an online report by NRE (''around 27 cm'') - but actually there is GonadState, Gonadweightand Statge Mature 3-7, in the historical data (how to interpret these?). StageMature 3-7 is 50% at length 31cm in the data (aggregated)'
$ws.Rows.Item(2).RowHeight = 43.2

# --- Insert a new row 16 for the "commercial year" issue/note pair.
#     This pushes the existing rows 16-22 down to 17-23 intact, preserving
#     their original cell layout (including the label-only row and the
#     region/large-region header rows). ---
$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16,1).Value = 'Not clear how to assign calendar year to commercial year'
$ws.Cells.Item(16,2).Value = 'Currently this is assumed to occur mostly in the second half, ie Nov 1 - Sep 1, so 2022/23 would be assigned the year 2023. '
$ws.Rows.Item(16).RowHeight = 28.8

# --- Expand the region / large-region notes (now at rows 19-22 after the insert). ---
$ws.Cells.Item(19,2).Value = 'Derwent Estuary, Tasman, Frederick Henry/Norfolk Bay, South-eastern coast, D''entrecasteaux Channel, South, Northwest Bay, SECest, SEC'
$ws.Cells.Item(20,2).Value = 'Great Oyster Bay, Central-eastern coast, Eastern coast, Coles Bay, Georges Bay, EC'
$ws.Cells.Item(21,2).Value = 'North-western coast, King Island, rocky cape, NWC'
$ws.Cells.Item(22,2).Value = 'Tamar River, North-eastern coast, Flinders Island, Spring Bay, Flinders/Eastcoast, NC, EC, Deal island, Hogan group, NEC, FI'

# --- Append a new row 24 for the "unknown" region code. ---
$ws.Cells.Item(24,1).Value = 'unknown'
$ws.Cells.Item(24,2).Value = 'EAT, ECS, ET, SET, CBS, no sample'

# --- Widen column B to fit the longer notes. ---
$ws.Columns.Item(2).ColumnWidth = 82

# --- Match the saved selection. ---
$ws.Range("B19").Select()
